$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

$ws.Cells.Item($row, 1).Value2 = "Globo"
$ws.Cells.Item($row, 2).Value2 = "RJ TV 2"
$ws.Cells.Item($row, 3).Value2 = "Obras"
$ws.Cells.Item($row, 4).Value2 = "2025-04-06T19:40"
$ws.Cells.Item($row, 5).Value2 = "Neutro"
$ws.Cells.Item($row, 6).Value2 = "A chuva no Estado do Rio. Defesas civis de Campos e cidades do Norte Fluminense acompanham deslocamento da frente fria. Repórter *ao vivo*. Locais mais atingidos de Campos foram Farol e Baixa Grande, na Baixada Campista; Dores de Macabu. Subsecretario da Defesa Civil, Edison Pessanha, disse que as chuvas não trouxeram nenhum transtorno grave. "
